$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 74076310
$ws.Range("I106").Value = 33335360
$ws.Range("K106").Value = 33335360
$ws.Range("M106").Value = -33334729

$ws.Range("H125").Value = 4237.316
$ws.Range("I125").Value = 12766
$ws.Range("J125").Value = 3233.9412
$ws.Range("K125").Value = 114894
$ws.Range("L125").Value = 29105.4708
$ws.Range("M125").Value = -112434
$ws.Range("N125").Value = -34025.4708

$ws.Range("H132").Value = 10102564
$ws.Range("I132").Value = 1260.3334
$ws.Range("J132").Value = 37039372
$ws.Range("K132").Value = 3781.0002
$ws.Range("L132").Value = 111118116
$ws.Range("M132").Value = -1251.0002
$ws.Range("N132").Value = -111123176

$ws.Range("H135").Value = 1547.8085
$ws.Range("I135").Value = 1240.275
$ws.Range("J135").Value = 3305.1428
$ws.Range("K135").Value = 11162.475
$ws.Range("L135").Value = 29746.2852
$ws.Range("M135").Value = -8627.475
$ws.Range("N135").Value = -34816.2852

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5241.9653
$ws.Range("I32").Value = 4064.7612
$ws.Range("J32").Value = 9185.6
$ws.Range("K32").Value = 4064.7612
$ws.Range("L32").Value = 9185.6
$ws.Range("M32").Value = -3777.7612
$ws.Range("N32").Value = -9759.6

$ws.Range("H61").Value = 4183.971
$ws.Range("I61").Value = 4331.727
$ws.Range("J61").Value = 1746
$ws.Range("K61").Value = 4331.727
$ws.Range("L61").Value = 1746
$ws.Range("M61").Value = -4119.727
$ws.Range("N61").Value = -2170

$ws.Range("H74").Value = 14286972
$ws.Range("I74").Value = 1035.25
$ws.Range("K74").Value = 1035.25
$ws.Range("M74").Value = -161.25

$ws.Range("H77").Value = 14286972
$ws.Range("I77").Value = 1035.25
$ws.Range("K77").Value = 5176.25
$ws.Range("M77").Value = -808.25

$ws.Range("H122").Value = 1028571.4
$ws.Range("I122").Value = 1511446.4
$ws.Range("J122").Value = 2461.875
$ws.Range("K122").Value = 4534339.199999999
$ws.Range("L122").Value = 7385.625
$ws.Range("M122").Value = -4531889.199999999
$ws.Range("N122").Value = -12285.625

$ws.Range("H123").Value = 29713.5
$ws.Range("J123").Value = 29713.5
$ws.Range("L123").Value = 29713.5
$ws.Range("N123").Value = -39513.5

$ws.Range("H132").Value = 3234.0857
$ws.Range("I132").Value = 2483.261
$ws.Range("J132").Value = 4673.1665
$ws.Range("K132").Value = 7449.782999999999
$ws.Range("L132").Value = 14019.4995
$ws.Range("M132").Value = -4919.782999999999
$ws.Range("N132").Value = -19079.4995

$ws.Range("H136").Value = 4183.971
$ws.Range("I136").Value = 4331.727
$ws.Range("J136").Value = 1746
$ws.Range("K136").Value = 12995.181
$ws.Range("L136").Value = 5238
$ws.Range("M136").Value = -10445.181
$ws.Range("N136").Value = -10338

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1500.3182
$ws.Range("I94").Value = 1033.1333
$ws.Range("J94").Value = 2501.4285
$ws.Range("K94").Value = 1033.1333
$ws.Range("L94").Value = 2501.4285
$ws.Range("M94").Value = -582.1333
$ws.Range("N94").Value = -3403.4285

$ws.Range("H134").Value = 2996.8508
$ws.Range("I134").Value = 3332.0205
$ws.Range("J134").Value = 2084.4443
$ws.Range("K134").Value = 9996.0615
$ws.Range("L134").Value = 6253.3329
$ws.Range("M134").Value = -7461.0615
$ws.Range("N134").Value = -11323.3329

$ws.Range("H140").Value = 44704.215
$ws.Range("J140").Value = 44704.215
$ws.Range("L140").Value = 44704.215
$ws.Range("N140").Value = -55064.215

$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4526652.5
$ws.Range("I16").Value = 10990397
$ws.Range("J16").Value = 2031.3
$ws.Range("K16").Value = 10990397
$ws.Range("L16").Value = 2031.3
$ws.Range("M16").Value = -10990110
$ws.Range("N16").Value = -2605.3

$ws.Range("H31").Value = 5996467
$ws.Range("I31").Value = 1609.4468
$ws.Range("J31").Value = 13221039
$ws.Range("K31").Value = 1609.4468
$ws.Range("L31").Value = 13221039
$ws.Range("M31").Value = -1314.4468
$ws.Range("N31").Value = -13221629

$ws.Range("H34").Value = 5996467
$ws.Range("I34").Value = 1609.4468
$ws.Range("J34").Value = 13221039
$ws.Range("K34").Value = 1609.4468
$ws.Range("L34").Value = 13221039
$ws.Range("M34").Value = -1407.4468
$ws.Range("N34").Value = -13221443

$ws.Range("H58").Value = 2646533.8
$ws.Range("I58").Value = 3788447.5
$ws.Range("J58").Value = 2101.842
$ws.Range("K58").Value = 3788447.5
$ws.Range("L58").Value = 2101.842
$ws.Range("M58").Value = -3788244.5
$ws.Range("N58").Value = -2507.842

$ws.Range("H105").Value = 13890796
$ws.Range("I105").Value = 16668540
$ws.Range("J105").Value = 2077.5
$ws.Range("K105").Value = 16668540
$ws.Range("L105").Value = 2077.5
$ws.Range("M105").Value = -16666793
$ws.Range("N105").Value = -5571.5

$ws.Range("H113").Value = 4526652.5
$ws.Range("I113").Value = 10990397
$ws.Range("J113").Value = 2031.3
$ws.Range("K113").Value = 10990397
$ws.Range("L113").Value = 2031.3
$ws.Range("M113").Value = -10988227
$ws.Range("N113").Value = -6371.3

$ws.Range("H132").Value = 4083464.5
$ws.Range("I132").Value = 5264749
$ws.Range("J132").Value = 2662.4546
$ws.Range("K132").Value = 15794247
$ws.Range("L132").Value = 7987.3638
$ws.Range("M132").Value = -15791717
$ws.Range("N132").Value = -13047.3638

$ws.Range("H134").Value = 7094251.5
$ws.Range("I134").Value = 11496771
$ws.Range("K134").Value = 34490313
$ws.Range("M134").Value = -34487778

$ws.Range("H136").Value = 2646533.8
$ws.Range("I136").Value = 3788447.5
$ws.Range("J136").Value = 2101.842
$ws.Range("K136").Value = 11365342.5
$ws.Range("L136").Value = 6305.526
$ws.Range("M136").Value = -11362792.5
$ws.Range("N136").Value = -11405.526

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 48471788
$ws.Range("I122").Value = 62637732
$ws.Range("J122").Value = 8334934.5
$ws.Range("K122").Value = 187913196
$ws.Range("L122").Value = 25004803.5
$ws.Range("M122").Value = -187910746
$ws.Range("N122").Value = -25009703.5

$ws.Range("H123").Value = 22000.125
$ws.Range("J123").Value = 22000.125
$ws.Range("L123").Value = 22000.125
$ws.Range("N123").Value = -26900.125

$ws.Range("H126").Value = 5770.1377
$ws.Range("I126").Value = 11801.1
$ws.Range("J126").Value = 2595.9473
$ws.Range("K126").Value = 35403.3
$ws.Range("L126").Value = 7787.841899999999
$ws.Range("M126").Value = -32933.3
$ws.Range("N126").Value = -12727.8419

$ws.Range("H136").Value = 14735.066
$ws.Range("J136").Value = 14735.066
$ws.Range("L136").Value = 44205.198
$ws.Range("N136").Value = -49305.198

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 6271474.5
$ws.Range("I122").Value = 11918062
$ws.Range("J122").Value = 1431542.9
$ws.Range("K122").Value = 35754186
$ws.Range("L122").Value = 4294628.699999999
$ws.Range("M122").Value = -35751736
$ws.Range("N122").Value = -4299528.699999999

$ws.Range("H132").Value = 12407224
$ws.Range("I132").Value = 15242635
$ws.Range("J132").Value = 2300
$ws.Range("K132").Value = 45727905
$ws.Range("L132").Value = 6900
$ws.Range("M132").Value = -45725375
$ws.Range("N132").Value = -11960

$ws.Range("H133").Value = 40323
$ws.Range("J133").Value = 40323
$ws.Range("L133").Value = 40323
$ws.Range("N133").Value = -45383

$ws.Range("H136").Value = 5704.228
$ws.Range("I136").Value = 4174.256
$ws.Range("J136").Value = 10403.429
$ws.Range("K136").Value = 12522.768
$ws.Range("L136").Value = 31210.287
$ws.Range("M136").Value = -9972.768
$ws.Range("N136").Value = -36310.287

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 3328.8
$ws.Range("I96").Value = 3066.6667
$ws.Range("J96").Value = 3722
$ws.Range("K96").Value = 3066.6667
$ws.Range("L96").Value = 3722
$ws.Range("M96").Value = -1693.6667
$ws.Range("N96").Value = -6468

$ws.Range("H113").Value = 1562.8096
$ws.Range("I113").Value = 1534.2354
$ws.Range("J113").Value = 1684.25
$ws.Range("K113").Value = 4602.706200000001
$ws.Range("L113").Value = 5052.75
$ws.Range("M113").Value = -2432.706200000001
$ws.Range("N113").Value = -9392.75

$ws.Range("H122").Value = 1188.4595
$ws.Range("I122").Value = 1174.5555
$ws.Range("J122").Value = 1226
$ws.Range("K122").Value = 3523.6665
$ws.Range("L122").Value = 3678
$ws.Range("M122").Value = -1073.6665
$ws.Range("N122").Value = -8578

$ws.Range("H123").Value = 40426
$ws.Range("J123").Value = 40426
$ws.Range("L123").Value = 40426
$ws.Range("N123").Value = -50226

$ws.Range("H126").Value = 1500.5555
$ws.Range("I126").Value = 600
$ws.Range("J126").Value = 1613.125
$ws.Range("K126").Value = 1800
$ws.Range("L126").Value = 4839.375
$ws.Range("M126").Value = 670
$ws.Range("N126").Value = -9779.375

$ws.Range("H132").Value = 901.0769
$ws.Range("I132").Value = 673.5405
$ws.Range("J132").Value = 1462.3334
$ws.Range("K132").Value = 2020.6215
$ws.Range("L132").Value = 4387.0002
$ws.Range("M132").Value = 509.3785000000003
$ws.Range("N132").Value = -9447.0002
